# HKD_YCONBootstrapping.xlsx - "few more fixes for the HKD curves"
#
# 1. Rename sheet "HKD_YCRH_Swaps_ON" -> "OIS" (cascades automatically to the
#    sheet-scoped defined name "DiscountingCurve" and to every RateHelpers!E*
#    formula that references the sheet by name).
# 2. General Settings!J4  : bump the "Trigger" timestamp.
# 3. General Settings!J10 : refresh the cached XML-export-path helper value.
# 4. OIS!T14              : move the sheet's remembered selection there.
# 5. OIS!M3               : refresh the cached ObjectHandler error message.
# 6. RateHelpers!G2       : refresh the cached HKD OIS quote value.
# 7. Try to restore the workbook window to a normal (non-hidden) state.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the curve sheet -------------------------------------------
$oisSheet = $wb.Worksheets.Item("HKD_YCRH_Swaps_ON")
$oisSheet.Name = "OIS"

# --- 2/3. General Settings literal value refreshes ------------------------
$gs = $wb.Worksheets.Item("General Settings")
$gs.Range("J4").Value = 41607.530034722222
$gs.Range("J10").Value = "N:\QuantLibXL-1.3.0\Data2\XML\"

# --- 4. Move OIS!s remembered selection to T14 -----------------------------
# Excel can only change the selection on the sheet that is active, so we
# briefly activate OIS, select the new cell, then restore the workbook's
# original active sheet ("Selected", tab index 3) so activeTab / tabSelected
# are left exactly as they were.
$originalActive = $wb.ActiveSheet
$oisSheet.Activate()
$oisSheet.Range("T14").Select()
$originalActive.Activate()

# --- 5. OIS!M3 cached ObjectHandler error message --------------------------
$oisSheet.Range("M3").Value = "ohObjectSave - Invalid parent path : N:\QuantLibXL-1.3.0\Data2\XML\HKD_YCONRH_Swaps.xml"

# --- 6. RateHelpers!G2 cached quote value -----------------------------------
$rh = $wb.Worksheets.Item("RateHelpers")
$rh.Range("G2").Value = 0.0006857

# --- 7. Un-hide the workbook window (best effort) ---------------------------
foreach ($win in $wb.Windows) {
    $win.Visible = $true
}
$excel.Windows.Item(1).Visible = $true
